$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tipps")

# Final results / tipps that became known for the games in rows 43, 44, 46
# (row 45's game result (column D) is still unknown - only the tipps got entered)
# Written as one block Range write so the engine's dependency graph picks up
# all five columns (D..H) together in a single pass.
$arr = New-Object 'object[,]' 4,5
$arr[0,0] = "NYG 10-07"; $arr[0,1] = "NE  17-09"; $arr[0,2] = "NE  17-13"; $arr[0,3] = "NE  14-09"; $arr[0,4] = "NE  19-14"
$arr[1,0] = "KC  31-17"; $arr[1,1] = "KC  24-21"; $arr[1,2] = "KC  23-17"; $arr[1,3] = "KC  28-17"; $arr[1,4] = "KC  27-17"
$arr[2,0] = $null;       $arr[2,1] = "MIN 23-17"; $arr[2,2] = "MIN 27-20"; $arr[2,3] = "MIN 21-14"; $arr[2,4] = "MIN 23-20"
$arr[3,0] = "LAR 37-14"; $arr[3,1] = "LAR 23-17"; $arr[3,2] = "LAR 22-20"; $arr[3,3] = "LAR 21-14"; $arr[3,4] = "LAR 23-21"
$ws.Range("D43:H46").Value = $arr
